{"js": "// Highlight quantitative impact metrics (percentages, dollar amounts, large\n// comma-grouped numbers, multipliers) in bold + dark-slate color across the\n// specific achievement / experience bullet paragraphs that contain them.\n//\n// Each matching paragraph is rebuilt as a sequence of runs: plain-text runs\n// stay untouched, metric runs get bold + color \"#2C3E50\".\n\n// Matches: \u00b187%, 23%, 73.5%, $400M, $4.7M, $840K, $1B, 1,200 (comma-grouped\n// integers), 2x / 1.5x multipliers. Evaluated in order, longest/most-specific\n// alternative first so e.g. \"$4.7M\" isn't split up.\nconst METRIC_RE = /\u00b1?\\$?\\d[\\d,]*\\.?\\d*%|\\$\\d[\\d,]*\\.?\\d*[MBK]|\\b\\d{1,3}(?:,\\d{3})+\\b|\\b\\d+\\.?\\d*x\\b/g;\n\n// The exact bullet/line texts targeted by this edit (each is matched against\n// a paragraph's full text so we only touch the intended paragraphs).\nconst TARGET_TEXTS = [\n  \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n  \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \u00b14.2% to \u00b12.1%\",\n  \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n  \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+\",\n  \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M\",\n  \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n\n  if (!TARGET_TEXTS.includes(text)) {\n    continue;\n  }\n\n  // Split the paragraph text into alternating plain / metric segments.\n  const segments = [];\n  let lastIndex = 0;\n  let match;\n  METRIC_RE.lastIndex = 0;\n  while ((match = METRIC_RE.exec(text)) !== null) {\n    if (match.index > lastIndex) {\n      segments.push({ text: text.slice(lastIndex, match.index), metric: false });\n    }\n    segments.push({ text: match[0], metric: true });\n    lastIndex = match.index + match[0].length;\n  }\n  if (lastIndex < text.length) {\n    segments.push({ text: text.slice(lastIndex), metric: false });\n  }\n\n  // No metrics found (shouldn't happen for our target list) -- skip.\n  if (!segments.some((s) => s.metric)) {\n    continue;\n  }\n\n  // Rebuild the paragraph contents as a sequence of runs.\n  para.clear();\n  await context.sync();\n\n  for (const seg of segments) {\n    const range = para.insertText(seg.text, Word.InsertLocation.end);\n    if (seg.metric) {\n      range.font.bold = true;\n      range.font.color = \"#2C3E50\";\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# Highlight quantitative impact metrics (percentages, dollar amounts, large\n# comma-grouped numbers, multipliers) in bold + dark-slate color (\"#2C3E50\")\n# across the specific achievement / experience bullet paragraphs that\n# contain them.\n#\n# Approach: for each target paragraph (matched by exact text), find every\n# metric substring with a regex, then bold + color just that sub-range of\n# the document (via absolute character offsets) \u2014 Word automatically splits\n# the underlying run(s) so only the matched text gets the new formatting.\n\n$d = $word.ActiveDocument\n\n# Regex matching percentages (incl. \"\u00b1\" and decimals), dollar amounts with\n# M/B/K suffixes, comma-grouped integers (e.g. 1,200), and \"2x\" multipliers.\n$pattern = [regex]'\\xB1?\\$?\\d[\\d,]*\\.?\\d*%|\\$\\d[\\d,]*\\.?\\d*[MBK]|\\b\\d{1,3}(?:,\\d{3})+\\b|\\b\\d+\\.?\\d*x\\b'\n\n# The exact bullet/line texts targeted by this edit (each is matched against\n# a paragraph's full text so we only touch the intended paragraphs).\n$plusMinus = [char]0x00B1\n$targets = @(\n  \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%\",\n  (\"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from \" + $plusMinus + \"4.2% to \" + $plusMinus + \"2.1%\"),\n  \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors for research platform development\",\n  \"\u2022 Created comprehensive meta-analysis framework handling millions of survey responses that became the `$400M Polling Consortium Database at The Analyst Institute, now valued at `$1B+\",\n  \"\u2022 Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M\",\n  \"\u2022 Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%\"\n)\n\nforeach ($p in $d.Paragraphs) {\n    $full = $p.Range\n    $ptext = $full.Text.TrimEnd([char]13, [char]7)\n\n    if (-not ($targets -contains $ptext)) {\n        continue\n    }\n\n    $matches = $pattern.Matches($ptext)\n    foreach ($m in $matches) {\n        $s = $full.Start + $m.Index\n        $e = $s + $m.Length\n        $metricRange = $d.Range($s, $e)\n        $metricRange.Font.Bold = 1\n        $metricRange.Font.Color = \"2C3E50\"\n    }\n}\n"}
